$d = $word.ActiveDocument

# The last paragraph in the document (empty placeholder right before the
# section break) currently has pStyle=ListParagraph, ind left=1080 and
# sz/szCs=28. It needs to become a numbered ("numId 3") list item carrying
# the new analysis text (sz/szCs=32), and a brand-new empty paragraph
# (ind left=720, sz/szCs=32, no style/numbering) must follow it.

$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count)
$r = $target.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="es-CO"/></w:rPr><w:t>Recibimos una entrada mediante un monitor serial para saber que patrón se quiere imprimir dependiendo de los valores lanzados por el potenciómetro que se dividirá en cuatro partes para determinar el patrón con el rango en el que se encuentra actualmente</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="es-CO"/></w:rPr></w:pPr></w:p>'

[void]$r.InsertXML($xml)
Write-Output "Replaced trailing empty paragraph with numbered analysis item + new blank paragraph."
